# Fruta / hortaliza, semanal
#
# A new weekly price record was inserted into the "Repollo" sheet at row 320,
# pushing the previously-existing rows 320-354 down to 321-355 (dimension
# grows from A1:R354 to A1:R355).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 320; everything below (old rows 320-354)
# shifts down to 321-355.
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new weekly record.
$ws.Range("A320").Value = 7
$ws.Range("B320").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C320").Value = "Ñuble"
$ws.Range("D320").Value = 45013
$ws.Range("E320").Value = 16
$ws.Range("F320").Value = 100112006
$ws.Range("G320").Value = "Repollo"
$ws.Range("H320").Value = "Crespo record"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 200
$ws.Range("K320").Value = 1300
$ws.Range("L320").Value = 1300
$ws.Range("M320").Value = 1300
$ws.Range("N320").Value = "$/unidad"
$ws.Range("O320").Value = "Provincia de Diguillín"
$ws.Range("P320").Value = 1300
$ws.Range("Q320").Value = 1
$ws.Range("R320").Value = "Hortaliza"
